$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only Taxonsorteringsordning (B) changes
$ws.Range("B2").Value = 96720

# Row 3 and Row 4 swap most of their content (A, D, E, F, G, H, I, J, AC),
# while column B gets new distinct values for each row.

# New Row 3 (previously held the "Blåsippa" record; now holds the "Knärot" record)
$ws.Range("A3").Value = 112234730
$ws.Range("B3").Value = 96720
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("I3").Value = "60"
$ws.Range("J3").Value = "stjälkar/strån/skott"
$ws.Range("AC3").Value = "Mer än 60 ex."

# New Row 4 (previously held the "Knärot" record; now holds the "Blåsippa" record)
$ws.Range("A4").Value = 112234736
$ws.Range("B4").Value = 98961
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 222498
$ws.Range("F4").Value = "Blåsippa"
$ws.Range("G4").Value = "Hepatica nobilis"
$ws.Range("H4").Value = "Schreb."
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("AC4").Value = ""
